$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a rolling IPO-tracking table (20 data rows, rows 2-21).
# This week's push adds two new listings ("넥스트바이오메디컬" and "유라클"),
# which pushes the whole table down and drops the oldest entry
# ("이노스페이스") off the bottom to keep the row count constant.

# 1) Make room for the two new rows by inserting blank rows at 6 and 8.
#    (This shifts existing rows down and preserves their original cell
#    types/formatting, notably the text-typed "확정공모가" values.)
$ws.Rows(6).Insert()
$ws.Rows(8).Insert()

# After both inserts, the old row 12 entry ("넥스트바이오메디컬") has moved to
# row 14 - it's now a duplicate of the fresh row we're about to write at
# row 6, so remove it there.
$ws.Rows(14).Delete()

# The oldest entry ("이노스페이스") has, after the shifts above, ended up at
# row 22 - drop it so the table stays at 20 data rows.
$ws.Rows(22).Delete()

# 2) Fill in the two new rows with this week's data.
$ws.Range("A6").Value = "넥스트바이오메디컬"
$ws.Range("B6").Value = "2024.07.29~08.02"
$ws.Range("C6").Value = "24,000~29,000"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 24000
$ws.Range("F6").Value = "한국투자증권"

$ws.Range("A8").Value = "유라클"
$ws.Range("B8").Value = "2024.07.29~08.02"
$ws.Range("C8").Value = "18,000~21,000"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 13518
$ws.Range("F8").Value = "키움증권"
